# Sync attendance_reports: reorder "Recorded By" (column G) entries so that
# System-generated entries ("System" / "system") are listed before the
# human/user entries (e.g. email addresses), preserving the relative order
# within each group.
#
# This mirrors an upstream sync that moved the automatic "System" markers
# to the front of the comma-separated "Recorded By" list on the
# "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows in column G whose "Recorded By" value needs to be reordered
# (matches the rows touched by the upstream sync commit).
$targetRows = @(
    2,3,5,6,7,8,10,11,12,13,14,15,17,18,19,20,21,22,24,26,
    28,29,31,32,33,34,36,37,38,39,40,41,43,44,45,46,47,48,50,52,
    54,55,57,58,59,60,62,63,64,65,66,67,69,70,71,72,73,74,76,78,
    80,81,82,83,84,85,86,90,92,93,94,96,99,101,
    106,107,108,109,110,111,112,116,118,119,120,122,125,127,
    132,133,134,135,136,137,138,142,144,145,146,148,151,153
)

foreach ($r in $targetRows) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $current = [string]$cell.Text

    if ([string]::IsNullOrEmpty($current)) { continue }

    $parts = $current -split ','
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $trimmed) {
        if ($p.ToLower() -eq 'system') {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    $reordered = $systemParts + $otherParts
    $newValue = [string]::Join(', ', $reordered)

    if ($newValue -ne $current) {
        $cell.Value = $newValue
    }
}
